# "Fruta / hortaliza, semanal" weekly update:
# a new week's price report (2021-09-21) is inserted for the two existing
# "Femacal de La Calera" Alcachofa series (Argentina(o)/Primera and
# Española/Extra) that currently sit in rows 125-126. That pushes every
# row from 127 downward by two, so we:
#   1. Insert two blank rows at 127:128 (shifts old 127..179 -> 129..181)
#   2. Copy the (still unmodified) old data from rows 125:126 down into
#      the freshly inserted rows 127:128
#   3. Overwrite rows 125:126 with the new week's Fecha/price figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("127:128").Insert()

$ws.Range("A125:R126").Copy()
$ws.Range("A127").PasteSpecial()

# Row 125: Argentina(o) / Primera
$ws.Range("D125").Value = "2021-09-21"
$ws.Range("K125").Value = 10000
$ws.Range("L125").Value = 10500
$ws.Range("M125").Value = 10238
$ws.Range("P125").Value = 205

# Row 126: Española / Extra
$ws.Range("D126").Value = "2021-09-21"
$ws.Range("K126").Value = 11000
$ws.Range("L126").Value = 11500
$ws.Range("M126").Value = 11262
$ws.Range("P126").Value = 375
